$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.649.62"
$ws.Range("E2").Value = "  -0.09%  "
$ws.Range("D3").Value = "1.642.90"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'215.33"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  +1.36%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.33%  "
$ws.Range("D9").Value = "'0.0628"
$ws.Range("E9").Value = "  +0.94%  "
$ws.Range("D10").Value = "'19.30"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").Value = "'0.0842"
$ws.Range("E11").Value = "  -0.02%  "
$ws.Range("D12").Value = "1.871.40"
$ws.Range("E12").Value = "  +0.51%  "
$ws.Range("D13").Value = "1.671.41"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("E14").Value = "  +2.13%  "
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("D16").Value = "'65.47"
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "26.697.19"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("E18").Value = "  +0.44%  "
$ws.Range("D19").Value = "'217.57"
$ws.Range("E19").Value = "  +0.70%  "
$ws.Range("E20").Value = "  +0.11%  "
$ws.Range("E21").Value = "  +1.30%  "
$ws.Range("E22").Value = "  +2.48%  "
$ws.Range("D23").Value = "'9.50"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").Value = "  +13.78%  "
$ws.Range("D25").Value = "'145.39"
$ws.Range("E25").Value = "  -1.53%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("E28").Value = "  +4.67%  "
$ws.Range("D29").Value = "'15.77"
$ws.Range("E29").Value = "  +1.42%  "
$ws.Range("E30").Value = "  +2.57%  "
$ws.Range("E31").Value = "  +0.81%  "
$ws.Range("D32").Value = "'3.39"
$ws.Range("E32").Value = "  +2.40%  "
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("D34").Value = "1.278.61"
$ws.Range("E34").Value = "  +4.25%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("E36").Value = "  +4.85%  "
$ws.Range("E37").Value = "  +0.60%  "
$ws.Range("D38").Value = "'0.532"
$ws.Range("E38").Value = "  +6.14%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'0.816"
$ws.Range("E41").Value = "  +2.65%  "
$ws.Range("D42").Value = "'2.26"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  +2.24%  "
$ws.Range("D44").Value = "1.782.09"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("D45").Value = "'91.71"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").Value = "'59.84"
$ws.Range("E46").Value = "  +8.25%  "
$ws.Range("E47").Value = "  +1.63%  "
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'7.79"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("D50").Value = "'0.0971"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "'0.407"
$ws.Range("E51").Value = "  -0.47%  "
